$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels + fill colors ---
$ws.Range("A1").Value = "Activity"
$ws.Range("B1").Value = "Date start:"
$ws.Range("C1").Value = "Date end:"
$ws.Range("D1").Value = "Time start:"
$ws.Range("E1").Value = "Time end:"

$ws.Range("A1").Interior.Color = 65535
$ws.Range("B1:C1").Interior.Color = 5296274
$ws.Range("D1:E1").Interior.Color = 255

# --- Second table header (row 17) ---
$ws.Range("A17").Value = "Activity "
$ws.Range("B17").Value = "Assigned to:"
$ws.Range("C17").Value = "Grade:"

$ws.Range("A17").Interior.Color = 65535
$ws.Range("B17").Interior.Color = 5296274
$ws.Range("C17").Interior.Color = 255

# --- Second table data (rows 18-22) ---
$ws.Range("A18").Value = "Creating  of forms - Login, registrarion, and main menu"
$ws.Range("B18").Value = "Miralles Renato"
$ws.Range("C18").Value = 84

$ws.Range("A19").Value = "Designing of forms - Login, registrarion, and main menu"
$ws.Range("B19").Value = "John Earl Azucena"
$ws.Range("C19").Value = 83

$ws.Range("A20").Value = "Creating database, beginning of coding in login, registration and main menu"
$ws.Range("B20").Value = "Irwin Legayo"
$ws.Range("C20").Value = 86

$ws.Range("A21").Value = "Research and studying about our features, gantt chart"
$ws.Range("B21").Value = "Norhana Daksla"
$ws.Range("C21").Value = 84

$ws.Range("A22").Value = "Also research and studying about the features and also gantt chart"
$ws.Range("B22").Value = "Alyssa Nikolei Co"
$ws.Range("C22").Value = 82

# --- Selection cursor moves to C14 ---
[void]$ws.Range("C14").Select()
